# Update automatico via Actualizar 02-20-2021 12-19-24
#
# The "Fecha" (Date) column D holds three stacked groups of 14 rows each,
# one group per monitored service-check run. Each automated refresh pushes
# a brand-new timestamp into the newest group (rows 2-15) and shifts the
# previously-newest timestamps down into the next group (rows 16-29), and
# those in turn shift down into the oldest group (rows 30-43), discarding
# whatever used to be there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Oldest group (rows 30-43) takes on what used to be the rows 16-29 value.
$ws.Range("D30:D43").Value = 44247.47070253472

# Middle group (rows 16-29) takes on what used to be the rows 2-15 value.
$ws.Range("D16:D29").Value = 44247.49201875

# Newest group (rows 2-15) gets the freshly recorded check timestamp.
$ws.Range("D2:D15").Value = 44247.51330931722
